# Update "想去人数" (F column) counts for several rows across three sheets.
# Values derived from the target diff (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 1025
    4  = 5754
    5  = 535
    6  = 977
    7  = 1004
    8  = 836
    10 = 42
    11 = 602
    12 = 38
    15 = 1940
    16 = 1500
    17 = 1017
    19 = 204
    20 = 372
    21 = 598
    22 = 207
    26 = 3274
    29 = 100
    30 = 143
    31 = 44
    32 = 446
    34 = 46
    38 = 781
    39 = 101
    41 = 71
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    4 = 245
    5 = 3
    6 = 176
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 1025
    6  = 5754
    7  = 535
    8  = 977
    10 = 245
    11 = 1004
    12 = 836
    13 = 3
    14 = 176
    16 = 42
    17 = 602
    18 = 38
    22 = 1940
    23 = 1500
    24 = 1017
    25 = 204
    26 = 372
    28 = 598
    29 = 207
    31 = 3274
    34 = 100
    35 = 143
    36 = 44
    37 = 446
    39 = 46
    42 = 781
    43 = 101
    45 = 71
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
